# ============================================================
# Edit script: restructure PlayerPerformance workbook
#   1. Insert new "Player Info" sheet (first)
#   2. Rename D1/B1 headers MATCH_CARD_LINK -> MATCH_CODE and replace
#      URL values with bare match codes on "ODI Batting"/"ODI Bowling"
#   3. Remove stray empty B25 cell on "ODI Batting"
#   4. Append new "ODI Batting Extra" sheet (last) with extra stats
# ============================================================

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------
# Step A: perform all sheet insertions FIRST (before grabbing any
# long-lived worksheet handles -- this engine resolves previously
# captured worksheet variables "live" by position, so capturing a
# handle and inserting a sheet before it can silently repoint the
# old variable at the newly inserted sheet). We re-fetch everything
# by name afterwards instead.
# ------------------------------------------------------------

# 1. Add "Player Info" sheet BEFORE "ODI Batting" (current index 1)
$newSheet1 = $wb.Worksheets.Add($wb.Worksheets.Item(1))
$newSheet1.Name = "Player Info"

# 2. Add "ODI Batting Extra" sheet AFTER the current last sheet ("ODI Bowling")
$lastIdx = $wb.Worksheets.Count
$newSheet2 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($lastIdx))
$newSheet2.Name = "ODI Batting Extra"

# ------------------------------------------------------------
# Step B: re-fetch stable references to every sheet we need, BY NAME
# ------------------------------------------------------------
$wsPlayerInfo = $wb.Worksheets.Item("Player Info")
$ws1 = $wb.Worksheets.Item("ODI Batting")
$ws2 = $wb.Worksheets.Item("ODI Bowling")
$ws4 = $wb.Worksheets.Item("ODI Batting Extra")

# ------------------------------------------------------------
# 3. Populate "Player Info"
# ------------------------------------------------------------
$wsPlayerInfo.Cells.Item(1,1).Value = "ID"
$wsPlayerInfo.Cells.Item(1,2).Value = "NAME"
$wsPlayerInfo.Cells.Item(1,3).Value = "BATTING_HAND"
$wsPlayerInfo.Cells.Item(1,4).Value = "BOWL_STYLE"

$hdrPI = $wsPlayerInfo.Range("A1:D1")
$hdrPI.Font.Bold = $true
$hdrPI.Borders.LineStyle = 1
$hdrPI.HorizontalAlignment = -4108
$hdrPI.VerticalAlignment = -4160

$wsPlayerInfo.Cells.Item(2,1).NumberFormat = "@"
$wsPlayerInfo.Cells.Item(2,1).Value = "4230"
$wsPlayerInfo.Cells.Item(2,2).Value = "Sohaib Maqsood"
$wsPlayerInfo.Cells.Item(2,3).Value = "Right Handed"
$wsPlayerInfo.Cells.Item(2,4).Value = "Does Not Bowl | Unknown"

# ------------------------------------------------------------
# 4. "ODI Batting": rename header + replace link column with match code
# ------------------------------------------------------------
$ws1.Cells.Item(1,4).Value = "MATCH_CODE"

$ws1.Cells.Item(2,4).NumberFormat = "@"
$ws1.Cells.Item(2,4).Value = "3575"
$ws1.Cells.Item(3,4).NumberFormat = "@"
$ws1.Cells.Item(3,4).Value = "3577"
$ws1.Cells.Item(4,4).NumberFormat = "@"
$ws1.Cells.Item(4,4).Value = "3582"
$ws1.Cells.Item(5,4).NumberFormat = "@"
$ws1.Cells.Item(5,4).Value = "3584"
$ws1.Cells.Item(6,4).NumberFormat = "@"
$ws1.Cells.Item(6,4).Value = "3585"
$ws1.Cells.Item(7,4).NumberFormat = "@"
$ws1.Cells.Item(7,4).Value = "3589"
$ws1.Cells.Item(8,4).NumberFormat = "@"
$ws1.Cells.Item(8,4).Value = "3590"
$ws1.Cells.Item(9,4).NumberFormat = "@"
$ws1.Cells.Item(9,4).Value = "3591"
$ws1.Cells.Item(10,4).NumberFormat = "@"
$ws1.Cells.Item(10,4).Value = "3592"
$ws1.Cells.Item(11,4).NumberFormat = "@"
$ws1.Cells.Item(11,4).Value = "3594"
$ws1.Cells.Item(12,4).NumberFormat = "@"
$ws1.Cells.Item(12,4).Value = "3618"
$ws1.Cells.Item(13,4).NumberFormat = "@"
$ws1.Cells.Item(13,4).Value = "3620"
$ws1.Cells.Item(14,4).NumberFormat = "@"
$ws1.Cells.Item(14,4).Value = "3624"
$ws1.Cells.Item(15,4).NumberFormat = "@"
$ws1.Cells.Item(15,4).Value = "3627"
$ws1.Cells.Item(16,4).NumberFormat = "@"
$ws1.Cells.Item(16,4).Value = "3658"
$ws1.Cells.Item(17,4).NumberFormat = "@"
$ws1.Cells.Item(17,4).Value = "3662"
$ws1.Cells.Item(18,4).NumberFormat = "@"
$ws1.Cells.Item(18,4).Value = "3666"
$ws1.Cells.Item(19,4).NumberFormat = "@"
$ws1.Cells.Item(19,4).Value = "3681"
$ws1.Cells.Item(20,4).NumberFormat = "@"
$ws1.Cells.Item(20,4).Value = "3751"
$ws1.Cells.Item(21,4).NumberFormat = "@"
$ws1.Cells.Item(21,4).Value = "3757"
$ws1.Cells.Item(22,4).NumberFormat = "@"
$ws1.Cells.Item(22,4).Value = "3770"
$ws1.Cells.Item(23,4).NumberFormat = "@"
$ws1.Cells.Item(23,4).Value = "3772"
$ws1.Cells.Item(24,4).NumberFormat = "@"
$ws1.Cells.Item(24,4).Value = "3776"
$ws1.Cells.Item(25,4).NumberFormat = "@"
$ws1.Cells.Item(25,4).Value = "3789"
$ws1.Cells.Item(26,4).NumberFormat = "@"
$ws1.Cells.Item(26,4).Value = "3792"
$ws1.Cells.Item(27,4).NumberFormat = "@"
$ws1.Cells.Item(27,4).Value = "3879"
$ws1.Cells.Item(28,4).NumberFormat = "@"
$ws1.Cells.Item(28,4).Value = "4472"
$ws1.Cells.Item(29,4).NumberFormat = "@"
$ws1.Cells.Item(29,4).Value = "4473"
$ws1.Cells.Item(30,4).NumberFormat = "@"
$ws1.Cells.Item(30,4).Value = "4476"

# remove stray empty cell B25 (present but blank in source workbook)
$ws1.Cells.Item(25,2).ClearContents()

# ------------------------------------------------------------
# 5. "ODI Bowling": rename header + replace link column with match code
# ------------------------------------------------------------
$ws2.Cells.Item(1,2).Value = "MATCH_CODE"

$ws2.Cells.Item(2,2).NumberFormat = "@"
$ws2.Cells.Item(2,2).Value = "3585"
$ws2.Cells.Item(3,2).NumberFormat = "@"
$ws2.Cells.Item(3,2).Value = "3627"
$ws2.Cells.Item(4,2).NumberFormat = "@"
$ws2.Cells.Item(4,2).Value = "3757"
$ws2.Cells.Item(5,2).NumberFormat = "@"
$ws2.Cells.Item(5,2).Value = "3772"

# ------------------------------------------------------------
# 6. Populate "ODI Batting Extra"
# ------------------------------------------------------------
$ws4.Cells.Item(1,1).Value = "MATCH_CODE"
$ws4.Cells.Item(1,2).Value = "BATTING_POSITION"
$ws4.Cells.Item(1,3).Value = "NUM_4"
$ws4.Cells.Item(1,4).Value = "NUM_6"
$ws4.Cells.Item(1,5).Value = "PERCENT_RUNS_OF_TOTAL"
$ws4.Cells.Item(1,6).Value = "MAN_OF_MATCH"

$hdrExtra = $ws4.Range("A1:F1")
$hdrExtra.Font.Bold = $true
$hdrExtra.Borders.LineStyle = 1
$hdrExtra.HorizontalAlignment = -4108
$hdrExtra.VerticalAlignment = -4160

# row 2: MATCH_CODE 3594
$ws4.Cells.Item(2,1).NumberFormat = "@"
$ws4.Cells.Item(2,1).Value = "3594"
$ws4.Cells.Item(2,2).NumberFormat = "@"
$ws4.Cells.Item(2,2).Value = ""
$ws4.Cells.Item(2,3).NumberFormat = "@"
$ws4.Cells.Item(2,3).Value = ""
$ws4.Cells.Item(2,4).NumberFormat = "@"
$ws4.Cells.Item(2,4).Value = ""
$ws4.Cells.Item(2,5).NumberFormat = "@"
$ws4.Cells.Item(2,5).Value = ""
$ws4.Cells.Item(2,6).Value = "NO"

# row 3: MATCH_CODE 3618
$ws4.Cells.Item(3,1).NumberFormat = "@"
$ws4.Cells.Item(3,1).Value = "3618"
$ws4.Cells.Item(3,2).Value = 4
$ws4.Cells.Item(3,3).NumberFormat = "@"
$ws4.Cells.Item(3,3).Value = "1"
$ws4.Cells.Item(3,4).NumberFormat = "@"
$ws4.Cells.Item(3,4).Value = "0"
$ws4.Cells.Item(3,5).NumberFormat = "@"
$ws4.Cells.Item(3,5).Value = "5.99%"
$ws4.Cells.Item(3,6).Value = "NO"

# row 4: MATCH_CODE 3620
$ws4.Cells.Item(4,1).NumberFormat = "@"
$ws4.Cells.Item(4,1).Value = "3620"
$ws4.Cells.Item(4,2).NumberFormat = "@"
$ws4.Cells.Item(4,2).Value = ""
$ws4.Cells.Item(4,3).NumberFormat = "@"
$ws4.Cells.Item(4,3).Value = ""
$ws4.Cells.Item(4,4).NumberFormat = "@"
$ws4.Cells.Item(4,4).Value = ""
$ws4.Cells.Item(4,5).NumberFormat = "@"
$ws4.Cells.Item(4,5).Value = ""
$ws4.Cells.Item(4,6).Value = "NO"

# row 5: MATCH_CODE 3624
$ws4.Cells.Item(5,1).NumberFormat = "@"
$ws4.Cells.Item(5,1).Value = "3624"
$ws4.Cells.Item(5,2).NumberFormat = "@"
$ws4.Cells.Item(5,2).Value = ""
$ws4.Cells.Item(5,3).NumberFormat = "@"
$ws4.Cells.Item(5,3).Value = ""
$ws4.Cells.Item(5,4).NumberFormat = "@"
$ws4.Cells.Item(5,4).Value = ""
$ws4.Cells.Item(5,5).NumberFormat = "@"
$ws4.Cells.Item(5,5).Value = ""
$ws4.Cells.Item(5,6).Value = "NO"

# row 6: MATCH_CODE 3627
$ws4.Cells.Item(6,1).NumberFormat = "@"
$ws4.Cells.Item(6,1).Value = "3627"
$ws4.Cells.Item(6,2).Value = 4
$ws4.Cells.Item(6,3).NumberFormat = "@"
$ws4.Cells.Item(6,3).Value = "0"
$ws4.Cells.Item(6,4).NumberFormat = "@"
$ws4.Cells.Item(6,4).Value = "0"
$ws4.Cells.Item(6,5).NumberFormat = "@"
$ws4.Cells.Item(6,5).Value = "0.61%"
$ws4.Cells.Item(6,6).Value = "NO"

# row 7: MATCH_CODE 3658
$ws4.Cells.Item(7,1).NumberFormat = "@"
$ws4.Cells.Item(7,1).Value = "3658"
$ws4.Cells.Item(7,2).Value = 7
$ws4.Cells.Item(7,3).NumberFormat = "@"
$ws4.Cells.Item(7,3).Value = "9"
$ws4.Cells.Item(7,4).NumberFormat = "@"
$ws4.Cells.Item(7,4).Value = "0"
$ws4.Cells.Item(7,5).NumberFormat = "@"
$ws4.Cells.Item(7,5).Value = "32.13%"
$ws4.Cells.Item(7,6).Value = "YES"

# row 8: MATCH_CODE 3662
$ws4.Cells.Item(8,1).NumberFormat = "@"
$ws4.Cells.Item(8,1).Value = "3662"
$ws4.Cells.Item(8,2).Value = 7
$ws4.Cells.Item(8,3).NumberFormat = "@"
$ws4.Cells.Item(8,3).Value = "0"
$ws4.Cells.Item(8,4).NumberFormat = "@"
$ws4.Cells.Item(8,4).Value = "1"
$ws4.Cells.Item(8,5).NumberFormat = "@"
$ws4.Cells.Item(8,5).Value = "3.86%"
$ws4.Cells.Item(8,6).Value = "NO"

# row 9: MATCH_CODE 3666
$ws4.Cells.Item(9,1).NumberFormat = "@"
$ws4.Cells.Item(9,1).Value = "3666"
$ws4.Cells.Item(9,2).Value = 7
$ws4.Cells.Item(9,3).NumberFormat = "@"
$ws4.Cells.Item(9,3).Value = "1"
$ws4.Cells.Item(9,4).NumberFormat = "@"
$ws4.Cells.Item(9,4).Value = "0"
$ws4.Cells.Item(9,5).NumberFormat = "@"
$ws4.Cells.Item(9,5).Value = "6.86%"
$ws4.Cells.Item(9,6).Value = "NO"

# row 10: MATCH_CODE 3681
$ws4.Cells.Item(10,1).NumberFormat = "@"
$ws4.Cells.Item(10,1).Value = "3681"
$ws4.Cells.Item(10,2).Value = 5
$ws4.Cells.Item(10,3).NumberFormat = "@"
$ws4.Cells.Item(10,3).Value = "2"
$ws4.Cells.Item(10,4).NumberFormat = "@"
$ws4.Cells.Item(10,4).Value = "1"
$ws4.Cells.Item(10,5).NumberFormat = "@"
$ws4.Cells.Item(10,5).Value = "14.78%"
$ws4.Cells.Item(10,6).Value = "NO"

# row 11: MATCH_CODE 3751
$ws4.Cells.Item(11,1).NumberFormat = "@"
$ws4.Cells.Item(11,1).Value = "3751"
$ws4.Cells.Item(11,2).Value = 5
$ws4.Cells.Item(11,3).NumberFormat = "@"
$ws4.Cells.Item(11,3).Value = "0"
$ws4.Cells.Item(11,4).NumberFormat = "@"
$ws4.Cells.Item(11,4).Value = "0"
$ws4.Cells.Item(11,5).NumberFormat = "@"
$ws4.Cells.Item(11,5).Value = ""
$ws4.Cells.Item(11,6).Value = "NO"

# row 12: MATCH_CODE 3757
$ws4.Cells.Item(12,1).NumberFormat = "@"
$ws4.Cells.Item(12,1).Value = "3757"
$ws4.Cells.Item(12,2).Value = 6
$ws4.Cells.Item(12,3).NumberFormat = "@"
$ws4.Cells.Item(12,3).Value = "4"
$ws4.Cells.Item(12,4).NumberFormat = "@"
$ws4.Cells.Item(12,4).Value = "1"
$ws4.Cells.Item(12,5).NumberFormat = "@"
$ws4.Cells.Item(12,5).Value = "31.25%"
$ws4.Cells.Item(12,6).Value = "NO"

# row 13: MATCH_CODE 3770
$ws4.Cells.Item(13,1).NumberFormat = "@"
$ws4.Cells.Item(13,1).Value = "3770"
$ws4.Cells.Item(13,2).Value = 7
$ws4.Cells.Item(13,3).NumberFormat = "@"
$ws4.Cells.Item(13,3).Value = "2"
$ws4.Cells.Item(13,4).NumberFormat = "@"
$ws4.Cells.Item(13,4).Value = "0"
$ws4.Cells.Item(13,5).NumberFormat = "@"
$ws4.Cells.Item(13,5).Value = "8.94%"
$ws4.Cells.Item(13,6).Value = "NO"

# row 14: MATCH_CODE 3772
$ws4.Cells.Item(14,1).NumberFormat = "@"
$ws4.Cells.Item(14,1).Value = "3772"
$ws4.Cells.Item(14,2).NumberFormat = "@"
$ws4.Cells.Item(14,2).Value = ""
$ws4.Cells.Item(14,3).NumberFormat = "@"
$ws4.Cells.Item(14,3).Value = ""
$ws4.Cells.Item(14,4).NumberFormat = "@"
$ws4.Cells.Item(14,4).Value = ""
$ws4.Cells.Item(14,5).NumberFormat = "@"
$ws4.Cells.Item(14,5).Value = ""
$ws4.Cells.Item(14,6).Value = "NO"

# row 15: MATCH_CODE 3776
$ws4.Cells.Item(15,1).NumberFormat = "@"
$ws4.Cells.Item(15,1).Value = "3776"
$ws4.Cells.Item(15,2).Value = 5
$ws4.Cells.Item(15,3).NumberFormat = "@"
$ws4.Cells.Item(15,3).Value = "1"
$ws4.Cells.Item(15,4).NumberFormat = "@"
$ws4.Cells.Item(15,4).Value = "0"
$ws4.Cells.Item(15,5).NumberFormat = "@"
$ws4.Cells.Item(15,5).Value = "3.60%"
$ws4.Cells.Item(15,6).Value = "NO"

# row 16: MATCH_CODE 3789
$ws4.Cells.Item(16,1).NumberFormat = "@"
$ws4.Cells.Item(16,1).Value = "3789"
$ws4.Cells.Item(16,2).Value = 6
$ws4.Cells.Item(16,3).NumberFormat = "@"
$ws4.Cells.Item(16,3).Value = ""
$ws4.Cells.Item(16,4).NumberFormat = "@"
$ws4.Cells.Item(16,4).Value = ""
$ws4.Cells.Item(16,5).NumberFormat = "@"
$ws4.Cells.Item(16,5).Value = ""
$ws4.Cells.Item(16,6).Value = "NO"

# row 17: MATCH_CODE 3792
$ws4.Cells.Item(17,1).NumberFormat = "@"
$ws4.Cells.Item(17,1).Value = "3792"
$ws4.Cells.Item(17,2).NumberFormat = "@"
$ws4.Cells.Item(17,2).Value = ""
$ws4.Cells.Item(17,3).NumberFormat = "@"
$ws4.Cells.Item(17,3).Value = ""
$ws4.Cells.Item(17,4).NumberFormat = "@"
$ws4.Cells.Item(17,4).Value = ""
$ws4.Cells.Item(17,5).NumberFormat = "@"
$ws4.Cells.Item(17,5).Value = ""
$ws4.Cells.Item(17,6).Value = "NO"

# row 18: MATCH_CODE 3879
$ws4.Cells.Item(18,1).NumberFormat = "@"
$ws4.Cells.Item(18,1).Value = "3879"
$ws4.Cells.Item(18,2).NumberFormat = "@"
$ws4.Cells.Item(18,2).Value = ""
$ws4.Cells.Item(18,3).NumberFormat = "@"
$ws4.Cells.Item(18,3).Value = ""
$ws4.Cells.Item(18,4).NumberFormat = "@"
$ws4.Cells.Item(18,4).Value = ""
$ws4.Cells.Item(18,5).NumberFormat = "@"
$ws4.Cells.Item(18,5).Value = ""
$ws4.Cells.Item(18,6).Value = "NO"

# row 19: MATCH_CODE 4472
$ws4.Cells.Item(19,1).NumberFormat = "@"
$ws4.Cells.Item(19,1).Value = "4472"
$ws4.Cells.Item(19,2).Value = 6
$ws4.Cells.Item(19,3).NumberFormat = "@"
$ws4.Cells.Item(19,3).Value = "1"
$ws4.Cells.Item(19,4).NumberFormat = "@"
$ws4.Cells.Item(19,4).Value = "1"
$ws4.Cells.Item(19,5).NumberFormat = "@"
$ws4.Cells.Item(19,5).Value = "13.48%"
$ws4.Cells.Item(19,6).Value = "NO"

# row 20: MATCH_CODE 4473
$ws4.Cells.Item(20,1).NumberFormat = "@"
$ws4.Cells.Item(20,1).Value = "4473"
$ws4.Cells.Item(20,2).Value = 6
$ws4.Cells.Item(20,3).NumberFormat = "@"
$ws4.Cells.Item(20,3).Value = "0"
$ws4.Cells.Item(20,4).NumberFormat = "@"
$ws4.Cells.Item(20,4).Value = "2"
$ws4.Cells.Item(20,5).NumberFormat = "@"
$ws4.Cells.Item(20,5).Value = "9.74%"
$ws4.Cells.Item(20,6).Value = "NO"

# row 21: MATCH_CODE 4476
$ws4.Cells.Item(21,1).NumberFormat = "@"
$ws4.Cells.Item(21,1).Value = "4476"
$ws4.Cells.Item(21,2).Value = 5
$ws4.Cells.Item(21,3).NumberFormat = "@"
$ws4.Cells.Item(21,3).Value = "2"
$ws4.Cells.Item(21,4).NumberFormat = "@"
$ws4.Cells.Item(21,4).Value = "0"
$ws4.Cells.Item(21,5).NumberFormat = "@"
$ws4.Cells.Item(21,5).Value = "2.42%"
$ws4.Cells.Item(21,6).Value = "NO"

Write-Host "Edit complete."
